$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, border, centered) from the existing header H1
# onto the two new header cells so they match the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new I and J numeric columns for data rows 2-19
$iVals = @(5, 6, 4, 8, 6, 7, 7, 5, 4, 1, 1, 1, 1, 1, 1, 1, 1, 7)
$jVals = @(6, 6, 5, 9, 7, 8, 7, 7, 6, 6, 5, 6, 5, 5, 6, 5, 6, 7)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}

Write-Output "I0 and IF columns added"
